$d = $word.ActiveDocument

$replacements = @(
    @{old="177×9=1593"; new="477×9=4293"},
    @{old="746×5=3730"; new="472×7=3304"},
    @{old="388×3=1164"; new="591×7=4137"},
    @{old="928×2=1856"; new="987×8=7896"},
    @{old="547×3=1641"; new="652×4=2608"},
    @{old="922×2=1844"; new="511×7=3577"},
    @{old="573×8=4584"; new="354×7=2478"},
    @{old="453×4=1812"; new="252×7=1764"},
    @{old="963×8=7704"; new="959×3=2877"},
    @{old="373×2=746";  new="148×8=1184"},
    @{old="219×4=876";  new="641×8=5128"},
    @{old="413×6=2478"; new="232×8=1856"},
    @{old="171×2=342";  new="322×4=1288"},
    @{old="629×2=1258"; new="689×2=1378"},
    @{old="178×8=1424"; new="498×9=4482"},
    @{old="618×4=2472"; new="889×5=4445"},
    @{old="831×9=7479"; new="212×6=1272"},
    @{old="762×4=3048"; new="626×7=4382"},
    @{old="956×9=8604"; new="222×4=888"},
    @{old="734×9=6606"; new="114×7=798"},
    @{old="247×2=494";  new="820×3=2460"},
    @{old="314×2=628";  new="596×6=3576"},
    @{old="511×2=1022"; new="694×6=4164"},
    @{old="651×7=4557"; new="284×8=2272"},
    @{old="575×8=4600"; new="637×2=1274"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
